$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(12, 1).Value = "2025-04-28 11:51:38"
$ws.Cells.Item(12, 2).Value = 236
